# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the
# ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with newly-fetched values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 10002510
$ws.Range("J45").Value = 10002510
$ws.Range("L45").Value = 30007530
$ws.Range("N45").Value = -30007914

$ws.Range("H70").Value = 2170.8
$ws.Range("I70").Value = 2170.8
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6512.400000000001
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -6242.400000000001
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 2170.8
$ws.Range("I73").Value = 2170.8
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6512.400000000001
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5576.400000000001
$ws.Range("N73").ClearContents()

$ws.Range("H98").Value = 1622.1034
$ws.Range("I98").Value = 1371.8889
$ws.Range("J98").Value = 5000
$ws.Range("K98").Value = 1371.8889
$ws.Range("L98").Value = 5000
$ws.Range("M98").Value = 126.1111000000001
$ws.Range("N98").Value = -7996

$ws.Range("H122").Value = 1622.1034
$ws.Range("I122").Value = 1371.8889
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4115.6667
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -1665.6667
$ws.Range("N122").Value = -19900

$ws.Range("H129").Value = 968.0909
$ws.Range("I129").Value = 636
$ws.Range("J129").Value = 1030.919
$ws.Range("K129").Value = 1908
$ws.Range("L129").Value = 3092.757000000001
$ws.Range("M129").Value = 3092
$ws.Range("N129").Value = -13092.757

$ws.Range("H138").Value = 5183.3706
$ws.Range("I138").Value = 1208.7028
$ws.Range("J138").Value = 8525.704
$ws.Range("K138").Value = 3626.1084
$ws.Range("L138").Value = 25577.112
$ws.Range("M138").Value = 1513.8916
$ws.Range("N138").Value = -35857.112

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 42926.25
$ws.Range("J121").Value = 42926.25
$ws.Range("L121").Value = 42926.25
$ws.Range("N121").Value = -46420.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 57152.89
$ws.Range("I20").Value = 60456
$ws.Range("J20").Value = 1000
$ws.Range("K20").Value = 60456
$ws.Range("L20").Value = 1000
$ws.Range("M20").Value = -60209
$ws.Range("N20").Value = -1494

$ws.Range("H22").Value = 290.2
$ws.Range("I22").Value = 183.66667
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 183.66667
$ws.Range("L22").Value = 450
$ws.Range("M22").Value = -10.66667000000001
$ws.Range("N22").Value = -796

$ws.Range("H56").Value = 17000
$ws.Range("J56").Value = 17000
$ws.Range("L56").Value = 17000
$ws.Range("N56").Value = -18478

$ws.Range("H86").Value = 35698.426
$ws.Range("I86").Value = 54140.047
$ws.Range("J86").Value = 3425.5833
$ws.Range("K86").Value = 54140.047
$ws.Range("L86").Value = 3425.5833
$ws.Range("M86").Value = -53017.047
$ws.Range("N86").Value = -5671.5833

$ws.Range("H89").Value = 35698.426
$ws.Range("I89").Value = 54140.047
$ws.Range("J89").Value = 3425.5833
$ws.Range("K89").Value = 270700.235
$ws.Range("L89").Value = 17127.9165
$ws.Range("M89").Value = -265084.235
$ws.Range("N89").Value = -28359.9165

$ws.Range("H94").Value = 843.7778
$ws.Range("I94").Value = 698
$ws.Range("J94").Value = 916.6667
$ws.Range("K94").Value = 698
$ws.Range("L94").Value = 916.6667
$ws.Range("M94").Value = -247
$ws.Range("N94").Value = -1818.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38718.652
$ws.Range("I31").Value = 1426.5
$ws.Range("J31").Value = 47242.57
$ws.Range("K31").Value = 1426.5
$ws.Range("L31").Value = 47242.57
$ws.Range("M31").Value = -1131.5
$ws.Range("N31").Value = -47832.57

$ws.Range("H34").Value = 38718.652
$ws.Range("I34").Value = 1426.5
$ws.Range("J34").Value = 47242.57
$ws.Range("K34").Value = 1426.5
$ws.Range("L34").Value = 47242.57
$ws.Range("M34").Value = -1224.5
$ws.Range("N34").Value = -47646.57

$ws.Range("H86").Value = 2506.7222
$ws.Range("I86").Value = 1962.375
$ws.Range("J86").Value = 2942.2
$ws.Range("K86").Value = 1962.375
$ws.Range("L86").Value = 2942.2
$ws.Range("M86").Value = -839.375
$ws.Range("N86").Value = -5188.2

$ws.Range("H89").Value = 2506.7222
$ws.Range("I89").Value = 1962.375
$ws.Range("J89").Value = 2942.2
$ws.Range("K89").Value = 9811.875
$ws.Range("L89").Value = 14711
$ws.Range("M89").Value = -4195.875
$ws.Range("N89").Value = -25943

$ws.Range("H132").Value = 42860730
$ws.Range("I132").Value = 37040612
$ws.Range("J132").Value = 62503624
$ws.Range("K132").Value = 111121836
$ws.Range("L132").Value = 187510872
$ws.Range("M132").Value = -111119306
$ws.Range("N132").Value = -187515932

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100000130
$ws.Range("I4").Value = 144.33333
$ws.Range("K4").Value = 432.99999
$ws.Range("M4").Value = -320.99999

$ws.Range("H131").Value = 864.25
$ws.Range("J131").Value = 869.5361
$ws.Range("L131").Value = 2608.6083
$ws.Range("N131").Value = -12688.6083

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 77316.17
$ws.Range("I70").Value = 149655.14
$ws.Range("J70").Value = 9799.799999999999
$ws.Range("K70").Value = 149655.14
$ws.Range("L70").Value = 9799.799999999999
$ws.Range("M70").Value = -149385.14
$ws.Range("N70").Value = -10339.8

$ws.Range("H73").Value = 77316.17
$ws.Range("I73").Value = 149655.14
$ws.Range("J73").Value = 9799.799999999999
$ws.Range("K73").Value = 149655.14
$ws.Range("L73").Value = 9799.799999999999
$ws.Range("M73").Value = -148719.14
$ws.Range("N73").Value = -11671.8

$ws.Range("H80").Value = 91004560
$ws.Range("I80").Value = 166840000
$ws.Range("J80").Value = 2030
$ws.Range("K80").Value = 166840000
$ws.Range("L80").Value = 2030
$ws.Range("M80").Value = -166839002
$ws.Range("N80").Value = -4026

$ws.Range("H83").Value = 91004560
$ws.Range("I83").Value = 166840000
$ws.Range("J83").Value = 2030
$ws.Range("K83").Value = 834200000
$ws.Range("L83").Value = 10150
$ws.Range("M83").Value = -834195008
$ws.Range("N83").Value = -20134

$ws.Range("H102").Value = 2819.1765
$ws.Range("I102").Value = 2683.4443
$ws.Range("J102").Value = 2971.875
$ws.Range("K102").Value = 2683.4443
$ws.Range("L102").Value = 2971.875
$ws.Range("M102").Value = -1061.4443
$ws.Range("N102").Value = -6215.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2961.5881
$ws.Range("I68").Value = 1733.3334
$ws.Range("J68").Value = 3224.7856
$ws.Range("K68").Value = 1733.3334
$ws.Range("L68").Value = 3224.7856
$ws.Range("M68").Value = -984.3334
$ws.Range("N68").Value = -4722.7856

$ws.Range("H71").Value = 2961.5881
$ws.Range("I71").Value = 1733.3334
$ws.Range("J71").Value = 3224.7856
$ws.Range("K71").Value = 8666.666999999999
$ws.Range("L71").Value = 16123.928
$ws.Range("M71").Value = -4922.666999999999
$ws.Range("N71").Value = -23611.928

$ws.Range("H122").Value = 2218
$ws.Range("I122").Value = 2080.923
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 6242.768999999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3792.768999999999
$ws.Range("N122").Value = -16900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13499.667
$ws.Range("J45").Value = 13499.667
$ws.Range("L45").Value = 13499.667
$ws.Range("N45").Value = -14481.667

$ws.Range("H107").Value = 84198.586
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 84198.586
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 252595.758
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -256435.758

$ws.Range("H122").Value = 1243
$ws.Range("I122").Value = 1158.5
$ws.Range("J122").Value = 1750
$ws.Range("K122").Value = 3475.5
$ws.Range("L122").Value = 5250
$ws.Range("M122").Value = -1025.5
$ws.Range("N122").Value = -10150

$ws.Range("H123").Value = 48500
$ws.Range("J123").Value = 48500
$ws.Range("L123").Value = 48500
$ws.Range("N123").Value = -58300
